$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "How many curves can I load in one go?"
$ws.Range("B10").Value = "openai"
$ws.Range("C10").Value = "You can load up to 450 curves at a time."

$ws.Range("A11").Value = "How many curves can I load in one go?"
$ws.Range("B11").Value = "llama3.2:latest"
$ws.Range("C11").Value = "You can load up to 450 curves at a time."
